# Bug fix: combo box matching on prefix-string.
# The task "if category A is removed, and category B is added, and category A
# is a starting-sub-string of B then all category A tasks are auto reset to B
# - they should not be" (Id 5) has been fixed, so it moves from the "Active"
# sheet (where it was Todo/Bug) to the "Inactive" sheet as a completed
# (Status=Done) task, with a "Done" date recorded.

$wb = $excel.ActiveWorkbook

$active = $wb.Worksheets.Item("Active")
$inactive = $wb.Worksheets.Item("Inactive")

# Remove the fixed task from the Active list (row 2 - Id 5).
$active.Rows.Item(2).Delete()

# Re-add it to the Inactive list as a completed task (new row 2), pushing the
# existing inactive rows down by one.
$inactive.Rows.Item(2).Insert()

$inactive.Range("A2").Value = 5
$inactive.Range("B2").Value = "if category A is removed, and category B is added," + [char]10 + "and category A is a starting-sub-string of B" + [char]10 + "then all category A tasks are auto reset to B - they should not be"

# Force these as text (matching the rest of the sheet) instead of letting
# Excel auto-convert the date-like strings into date serials.
$inactive.Range("C2:F2").NumberFormat = "@"
$inactive.Range("C2").Value = "Done"
$inactive.Range("D2").Value = "Bug"
$inactive.Range("E2").Value = "12/1/2017"
$inactive.Range("F2").Value = "4/11/2018"

# Restore default styling on the new row.
$inactive.Range("A2:F2").Style = "Normal"
